$d = $word.ActiveDocument

# Fix typo: "This differences" -> "These differences"
# (Find.Execute args: FindText, MatchCase, MatchWholeWord, MatchWildcards,
#  MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)
$null = $d.Content.Find.Execute("This", $true, $true, $false, $false, $false,
                                 $true, 1, $false, "These", 2)
